$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '60.091.93'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +2.66%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '3.204.00'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +1.50%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '536.96'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -0.02%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '145.13'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +3.74%  '

$ws.Range("D7").NumberFormat = '@'
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.529'
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  +2.96%  '

$ws.Range("E9").Value = '  -0.23%  '

$ws.Range("E10").Value = '  +2.89%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.432'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +2.74%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '3.755.82'
$ws.Range("D12").Style = 'Normal'
$ws.Range("E12").Value = '  +1.47%  '

$ws.Range("E13").Value = '  -1.22%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '25.91'
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  +0.50%  '

$ws.Range("E15").Value = '  +1.69%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '60.144.01'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +2.64%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '3.207.18'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +1.52%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '6.25'
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  +0.31%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '13.16'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +0.94%  '

$ws.Range("E20").Value = '  +0.54%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '375.67'
$ws.Range("D21").Style = 'Normal'
$ws.Range("E21").Value = '  +0.71%  '

$ws.Range("E22").Value = '  +0.31%  '

$ws.Range("E23").Value = '  +1.91%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '70.03'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("E25").Value = '  +1.26%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '8.77'
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  +8.85%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +0.46%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '0.0₃0897'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +2.99%  '

$ws.Range("E29").Value = '  +0.65%  '

$ws.Range("E30").Value = '  +1.56%  '

$ws.Range("E31").Value = '  -0.52%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '5.38'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +3.61%  '

$ws.Range("E33").Value = '  +7.40%  '

$ws.Range("E34").Value = '  +2.45%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '156.60'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -2.16%  '

$ws.Range("E36").Value = '  -1.66%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '2.804.51'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +6.32%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '25.64'
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  +1.38%  '

$ws.Range("E39").Value = '  +3.16%  '

$ws.Range("E40").Value = '  +0.73%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '4.23'
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +1.53%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '39.88'
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  +2.84%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '0.0294'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +4.66%  '

$ws.Range("E44").Value = '  +1.16%  '

$ws.Range("E45").Value = '  +3.37%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '3.246.43'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  +1.48%  '

$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '0.985'
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  +0.21%  '

$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '0.812'
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  +7.00%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '6.16'
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -0.96%  '

$ws.Range("E51").Value = '  +0.02%  '
